$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading "'" forces the numeric-looking value to stay stored as text,
# matching the source data's inline-string "Price" column.

# --- Price (column D) refresh for rows whose coin identity is unchanged ---
$ws.Range("D2").Value  = "'247.29"
$ws.Range("D3").Value  = "'26.44"
$ws.Range("D4").Value  = "'5.080"
$ws.Range("D5").Value  = "'0.05615"
$ws.Range("D7").Value  = "'0.8133"
$ws.Range("D8").Value  = "'0.8488"
$ws.Range("D19").Value = "'0.06995"
$ws.Range("D20").Value = "'0.03176"
$ws.Range("D22").Value = "'3.741"
$ws.Range("D23").Value = "'0.04659"
$ws.Range("D24").Value = "'0.1350"
$ws.Range("D26").Value = "'0.004613"
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("D43").Value = "'0.002599"
$ws.Range("D45").Value = "'0.00005296"
$ws.Range("D48").Value = "'0.002568"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"

# --- Rows 9-18: coin-ranking reshuffle (each coin moves up one slot, WazirX wraps to the bottom) ---
$ws.Range("B9").Value  = "BitrueCoin"
$ws.Range("C9").Value  = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value  = "'0.02834"
$ws.Range("E9").Value  = "8BitrueCoinBTR"

$ws.Range("B10").Value = "BitMartToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "'0.09396"
$ws.Range("E10").Value = "9BitMartTokenBMX"

$ws.Range("B11").Value = "BitForexToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "'0.001516"
$ws.Range("E11").Value = "10BitForexTokenBF"

$ws.Range("B12").Value = "One"
$ws.Range("C12").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "'0.0006000"
$ws.Range("E12").Value = "11OneONE"

$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006222"
$ws.Range("E13").Value = "12TigerCashTCH"

$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.586"
$ws.Range("E14").Value = "13LEOLEO"

$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.058"
$ws.Range("E15").Value = "14GateTokenGT"

$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.118"
$ws.Range("E16").Value = "15BTSETokenBTSE"

$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3181"
$ws.Range("E17").Value = "16BitpandaEcosystemTokenBEST"

$ws.Range("B18").Value = "WazirX"
$ws.Range("C18").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D18").Value = "'0.1349"
$ws.Range("E18").Value = "17WazirXWRX"

# --- Price + 24h extremum-label updates further down the table ---
$ws.Range("D41").Value = "'0.006120"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

$ws.Range("D44").Value = "'0.008575"
$ws.Range("E44").Value = "43LocalTradersLCT"

$ws.Range("D47").Value = "'0.1200"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
